# Auto-generated edit script applying value updates to Kraken_Profits sheets
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 636.125
$ws.Range("I19").Value = 756.75
$ws.Range("K19").Value = 756.75
$ws.Range("M19").Value = -581.75
$ws.Range("H32").Value = 7109.2
$ws.Range("J32").Value = 7164.8
$ws.Range("L32").Value = 7164.8
$ws.Range("N32").Value = -7816.8
$ws.Range("H64").Value = 3083.6667
$ws.Range("J64").Value = 4004
$ws.Range("L64").Value = 4004
$ws.Range("N64").Value = -4500
$ws.Range("H67").Value = 3083.6667
$ws.Range("J67").Value = 4004
$ws.Range("L67").Value = 4004
$ws.Range("N67").Value = -5720
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4064
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -20320
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184
$ws.Range("H88").Value = 3274.6365
$ws.Range("I88").Value = 3849.8333
$ws.Range("J88").Value = 2584.4
$ws.Range("K88").Value = 3849.8333
$ws.Range("L88").Value = 2584.4
$ws.Range("M88").Value = -3443.8333
$ws.Range("N88").Value = -3396.4
$ws.Range("H91").Value = 3274.6365
$ws.Range("I91").Value = 3849.8333
$ws.Range("J91").Value = 2584.4
$ws.Range("K91").Value = 3849.8333
$ws.Range("L91").Value = 2584.4
$ws.Range("M91").Value = -2445.8333
$ws.Range("N91").Value = -5392.4
$ws.Range("H100").Value = 3587.25
$ws.Range("I100").Value = 2233
$ws.Range("J100").Value = 4399.8
$ws.Range("K100").Value = 2233
$ws.Range("L100").Value = 4399.8
$ws.Range("M100").Value = -1692
$ws.Range("N100").Value = -5481.8
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 12500000
$ws.Range("I13").Value = 12500000
$ws.Range("K13").Value = 12500000
$ws.Range("M13").Value = -12499856
$ws.Range("H37").Value = 12750
$ws.Range("I37").Value = 5500
$ws.Range("J37").Value = 20000
$ws.Range("K37").Value = 5500
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = -5227
$ws.Range("N37").Value = -20546
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 5000
$ws.Range("M63").Value = -4314
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 25000
$ws.Range("M66").Value = -21568
$ws.Range("H74").Value = 5169
$ws.Range("I74").Value = 5169
$ws.Range("K74").Value = 5169
$ws.Range("M74").Value = -4295
$ws.Range("H77").Value = 5169
$ws.Range("I77").Value = 5169
$ws.Range("K77").Value = 25845
$ws.Range("M77").Value = -21477

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H99").Value = 4285.4287
$ws.Range("I99").Value = 4285.4287
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4285.4287
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2787.4287
$ws.Range("N99").ClearContents()

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2493.6667
$ws.Range("I99").Value = 2493.6667
$ws.Range("K99").Value = 2493.6667
$ws.Range("M99").Value = -995.6667000000002
$ws.Range("H126").Value = 2493.6667
$ws.Range("I126").Value = 2493.6667
$ws.Range("K126").Value = 7481.000100000001
$ws.Range("M126").Value = -5011.000100000001

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 518
$ws.Range("J107").Value = 749.5
$ws.Range("L107").Value = 2248.5
$ws.Range("N107").Value = -6088.5
$ws.Range("H113").Value = 1098.25
$ws.Range("J113").Value = 1098.25
$ws.Range("L113").Value = 3294.75
$ws.Range("N113").Value = -7634.75
$ws.Range("H131").Value = 3000
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080
$ws.Range("H132").Value = 1126.1818
$ws.Range("I132").Value = 1126.1818
$ws.Range("K132").Value = 10135.6362
$ws.Range("M132").Value = -7605.636200000001

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51498
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -157488
$ws.Range("H140").Value = 79999
$ws.Range("J140").Value = 79999
$ws.Range("L140").Value = 79999
$ws.Range("N140").Value = -90359

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H16").Value = 1723.375
$ws.Range("I16").Value = 1881.5
$ws.Range("J16").Value = 1249
$ws.Range("K16").Value = 1881.5
$ws.Range("L16").Value = 1249
$ws.Range("M16").Value = -1589
$ws.Range("H40").Value = 3389
$ws.Range("I40").Value = 3389
$ws.Range("K40").Value = 3389
$ws.Range("M40").Value = -3253
$ws.Range("H122").Value = 5954.6665
$ws.Range("I122").Value = 5954.6665
$ws.Range("K122").Value = 17863.9995
$ws.Range("M122").Value = -15413.9995

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 404.33334
$ws.Range("I113").Value = 329.875
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 989.625
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1180.375
$ws.Range("N113").Value = -7340

Write-Host "Applied all updates to Kraken_Profits sheets"
